# Update "想去人数" (interest count) values across the relevant worksheets.
# Sheet "展览" (Worksheets index/name)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1115
$ws1.Range("F5").Value = 2738
$ws1.Range("F7").Value = 690
$ws1.Range("F11").Value = 694
$ws1.Range("F12").Value = 100
$ws1.Range("F14").Value = 1594

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2019
$ws3.Range("F5").Value = 252

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2019
$ws4.Range("F5").Value = 252
$ws4.Range("F12").Value = 1115
$ws4.Range("F16").Value = 2738
$ws4.Range("F22").Value = 690
$ws4.Range("F27").Value = 694
$ws4.Range("F28").Value = 100
$ws4.Range("F31").Value = 1594
